# Updates cryptos list figures (Price column D, Volume(1h) column E)
# on Sheet1, matching the scraped refresh committed by the GitHub Actions
# job on Fri Aug  9 03:47:38 UTC 2024.
#
# Price column (D) cells are stored as text in the source data (e.g.
# "60.768.30", "0.994", "513.54"), so NumberFormat is forced to "@" (Text)
# before each assignment to stop Excel from auto-coercing the string into
# a number (which would silently drop values like the trailing zero in
# "60.50" -> 60.5, or re-interpret "2.639.49"-style thousand-grouped
# figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.768.30"
$ws.Range("E2").Value = "  +6.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.638.13"
$ws.Range("E3").Value = "  +7.90%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.54"
$ws.Range("E5").Value = "  +5.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.84"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.994"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.680.21"
$ws.Range("E9").Value = "  +9.83%  "
$ws.Range("E10").Value = "  +8.93%  "
$ws.Range("E11").Value = "  +6.20%  "
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.112.93"
$ws.Range("E14").Value = "  +9.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.882.19"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.24"
$ws.Range("E16").Value = "  +6.95%  "
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.680.65"
$ws.Range("E18").Value = "  +9.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.83"
$ws.Range("E20").Value = "  +6.26%  "
$ws.Range("E21").Value = "  +5.87%  "
$ws.Range("E22").Value = "  +4.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.50"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.788.21"
$ws.Range("E26").Value = "  +9.99%  "
$ws.Range("E27").Value = "  +5.31%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0870"
$ws.Range("E29").Value = "  +10.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.58"
$ws.Range("E30").Value = "  +4.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.67"
$ws.Range("E32").Value = "  +4.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.99"
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.74"
$ws.Range("E35").Value = "  +8.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.09"
$ws.Range("E36").Value = "  +10.47%  "
$ws.Range("E37").Value = "  +6.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("E39").Value = "  +12.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "309.15"
$ws.Range("E40").Value = "  +14.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.79"
$ws.Range("E41").Value = "  +7.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.840"
$ws.Range("E42").Value = "  +30.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.71"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("E44").Value = "  +8.87%  "
$ws.Range("E45").Value = "  +8.71%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.33"
$ws.Range("E47").Value = "  +15.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.03"
$ws.Range("E48").Value = "  +7.16%  "
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0237"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.039.01"
$ws.Range("E51").Value = "  +9.77%  "

